# Apply the updated symbol-list values (price, 1h volume %, and hour)
# scraped for this run. Values are written as literal text (matching the
# source sheet's inlineStr cells), not as numbers/percentages, so a leading
# apostrophe forces text entry; Style is reset to Normal afterwards so no
# stray "quote prefix" number format lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{Cell="D2"; Value="288.72"},
    @{Cell="E2"; Value="0.11%"},
    @{Cell="G2"; Value="14"},
    @{Cell="D3"; Value="31.04"},
    @{Cell="E3"; Value="2.13%"},
    @{Cell="G3"; Value="14"},
    @{Cell="D4"; Value="4.962"},
    @{Cell="E4"; Value="0.34%"},
    @{Cell="G4"; Value="14"},
    @{Cell="D5"; Value="0.07348"},
    @{Cell="E5"; Value="1.24%"},
    @{Cell="G5"; Value="14"},
    @{Cell="D6"; Value="2.364"},
    @{Cell="E6"; Value="31.55%"},
    @{Cell="G6"; Value="14"},
    @{Cell="D7"; Value="7.730"},
    @{Cell="E7"; Value="1.82%"},
    @{Cell="G7"; Value="14"},
    @{Cell="D8"; Value="0.9118"},
    @{Cell="E8"; Value="1.01%"},
    @{Cell="G8"; Value="14"},
    @{Cell="E9"; Value="19.11%"},
    @{Cell="G9"; Value="14"},
    @{Cell="D10"; Value="0.1705"},
    @{Cell="E10"; Value="2.34%"},
    @{Cell="G10"; Value="14"},
    @{Cell="E11"; Value="2.18%"},
    @{Cell="G11"; Value="14"},
    @{Cell="E12"; Value="1.89%"},
    @{Cell="G12"; Value="14"},
    @{Cell="D13"; Value="0.09970"},
    @{Cell="E13"; Value="-0.31%"},
    @{Cell="G13"; Value="14"},
    @{Cell="D14"; Value="0.001495"},
    @{Cell="E14"; Value="-0.12%"},
    @{Cell="G14"; Value="14"},
    @{Cell="D15"; Value="0.005754"},
    @{Cell="E15"; Value="-1.04%"},
    @{Cell="G15"; Value="14"},
    @{Cell="D16"; Value="3.471"},
    @{Cell="E16"; Value="-0.05%"},
    @{Cell="G16"; Value="14"},
    @{Cell="D17"; Value="3.729"},
    @{Cell="E17"; Value="0.63%"},
    @{Cell="G17"; Value="14"},
    @{Cell="D18"; Value="2.106"},
    @{Cell="E18"; Value="1.51%"},
    @{Cell="G18"; Value="14"},
    @{Cell="D19"; Value="0.3323"},
    @{Cell="E19"; Value="0.71%"},
    @{Cell="G19"; Value="14"},
    @{Cell="E20"; Value="-0.47%"},
    @{Cell="G20"; Value="14"},
    @{Cell="D21"; Value="4.170"},
    @{Cell="E21"; Value="5.23%"},
    @{Cell="G21"; Value="14"},
    @{Cell="D22"; Value="0.2101"},
    @{Cell="E22"; Value="0.12%"},
    @{Cell="G22"; Value="14"},
    @{Cell="D23"; Value="0.04528"},
    @{Cell="G23"; Value="14"},
    @{Cell="E24"; Value="-0.17%"},
    @{Cell="G24"; Value="14"},
    @{Cell="D25"; Value="0.004172"},
    @{Cell="E25"; Value="-10.03%"},
    @{Cell="G25"; Value="14"},
    @{Cell="D26"; Value="0.0001300"},
    @{Cell="E26"; Value="0.03%"},
    @{Cell="G26"; Value="14"},
    @{Cell="D27"; Value="0.0003395"},
    @{Cell="G27"; Value="14"},
    @{Cell="G28"; Value="14"},
    @{Cell="G29"; Value="14"},
    @{Cell="G30"; Value="14"},
    @{Cell="G31"; Value="14"},
    @{Cell="G32"; Value="14"},
    @{Cell="G33"; Value="14"},
    @{Cell="G34"; Value="14"},
    @{Cell="G35"; Value="14"},
    @{Cell="G36"; Value="14"},
    @{Cell="G37"; Value="14"},
    @{Cell="G38"; Value="14"},
    @{Cell="D39"; Value="0.01578"},
    @{Cell="E39"; Value="0.90%"},
    @{Cell="G39"; Value="14"},
    @{Cell="D40"; Value="0.04471"},
    @{Cell="G40"; Value="14"},
    @{Cell="D41"; Value="0.007377"},
    @{Cell="E41"; Value="0.94%"},
    @{Cell="G41"; Value="14"},
    @{Cell="D42"; Value="0.009861"},
    @{Cell="E42"; Value="-1.95%"},
    @{Cell="G42"; Value="14"},
    @{Cell="D43"; Value="0.1331"},
    @{Cell="E43"; Value="2.13%"},
    @{Cell="G43"; Value="14"},
    @{Cell="E44"; Value="11.43%"},
    @{Cell="G44"; Value="14"},
    @{Cell="D45"; Value="0.008778"},
    @{Cell="E45"; Value="-6.81%"},
    @{Cell="G45"; Value="14"},
    @{Cell="D46"; Value="0.00006112"},
    @{Cell="E46"; Value="3.88%"},
    @{Cell="G46"; Value="14"},
    @{Cell="E47"; Value="0.17%"},
    @{Cell="G47"; Value="14"},
    @{Cell="D48"; Value="2.606"},
    @{Cell="E48"; Value="15.58%"},
    @{Cell="G48"; Value="14"},
    @{Cell="D49"; Value="0.002001"},
    @{Cell="E49"; Value="-33.30%"},
    @{Cell="G49"; Value="14"},
    @{Cell="E50"; Value="0.17%"},
    @{Cell="G50"; Value="14"},
    @{Cell="E51"; Value="0.17%"},
    @{Cell="G51"; Value="14"}
)

foreach ($edit in $edits) {
    $cell = $ws.Range($edit.Cell)
    # Leading apostrophe forces text storage even for numeric-looking strings
    $cell.Value = "'" + $edit.Value
    # Drop the auto-applied quote-prefix style so formatting matches the source
    $cell.Style = "Normal"
}

